$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update status column to "OK" for the first three existing requirement rows ---
$ws.Range("C2").Value = "OK"
$ws.Range("C3").Value = "OK"
$ws.Range("C4").Value = "OK"

# --- Update the wording of the third requirement row (B4) ---
$ws.Range("B4").Value = "c#程序读取数据库，列出连续减少N期的股票`n单击股票后，右边列出每期的股东人数信息"

# --- Copy formatting from row 3 down to the new rows 5-7 ---
$ws.Range("A3:C3").Copy()
$ws.Range("A5:C7").PasteSpecial(-4122)

# --- Add three new requirement rows ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "输入号码和股票名可以定位股票"
$ws.Range("C5").Value = "OK"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "可以选择最后一期的日期"
$ws.Range("C6").Value = "OK"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "连续减少期数设为0，表示列出所有股票"
$ws.Range("C7").Value = "OK"

# --- Update the active selection to match the target state ---
$ws.Range("C8").Select() | Out-Null
